# Apply "Add configurable invalid path detection with detailed error messages" edit
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Export")

$newScanDate = "2025-09-10 21:53:58"

# Update the ScanDate (column J) for all existing data rows (2-8) to the new timestamp
for ($r = 2; $r -le 8; $r++) {
    $ws.Cells.Item($r, 10).Value = $newScanDate
}

# Row 5: path is a directory -> give a specific ScanError message
$ws.Cells.Item(5, 8).Value = "Path is a directory, not a file"

# Row 7: empty file path -> give a specific ScanError message
$ws.Cells.Item(7, 8).Value = "Empty file path"

# Row 8: path marked as N\A -> give a specific ScanError message
$ws.Cells.Item(8, 8).Value = "Path marked as N/A"

# Add new row 9 for the jackson-databind jar (JAR filename validation fix:
# a path under "C:\Program Files\..." no longer falsely flagged, so FileExists = N)
$ws.Cells.Item(9, 1).Value = "LPRIME"
$ws.Cells.Item(9, 2).Value = "Windows Server 2019"
$ws.Cells.Item(9, 3).Value = "C:\Program Files\BMC Software\Control-M Agent\Default\EXE_9.0.20.200\Jars\jackson-databind-2.10.3.jar"
$ws.Cells.Item(9, 4).Value = "X"
$ws.Cells.Item(9, 5).Value = "N"
$ws.Cells.Item(9, 6).Value = ""
$ws.Cells.Item(9, 7).Value = ""
$ws.Cells.Item(9, 8).Value = ""
$ws.Cells.Item(9, 10).Value = $newScanDate

# Adjust the used range / selection to include the new row
$ws.Range("E2:K9").Select()

# Widen column C to fit the new, longer path text (and drop the old bestFit autosize)
$ws.Columns.Item(3).ColumnWidth = 71.453125
